$p = $ppt.ActivePresentation

# 1. Add a Custom Show ("Custom Show 1") that contains the first slide.
$s1 = $p.Slides.Item(1)
$slideIds = @($s1.SlideID)
$p.SlideShowSettings.NamedSlideShows.Add("Custom Show 1", $slideIds)

# 2. Refresh the cached "datetimeFigureOut" date placeholder text
#    (slide master + every slide layout) from 3/18/2020 to 2/13/2021.
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "3/18/2020") {
                $sh.TextFrame.TextRange.Text = "2/13/2021"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder($master.Shapes)

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder($layout.Shapes)
}
